$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (item "9", AC-RFD-00-235-COP-H2B1-01) gets merged/renamed into the
#     new "AC-RFD-00-235-000-H3B1-01" model, its Kullanim Yeri / Baglanti Turu
#     change and it gains a Link hyperlink. ---
$ws.Range("B10").Value = "AC-RFD-00-235-000-H3B1-01"
$ws.Range("F10").Value = "Genel"
$ws.Range("H10").Value = "Kablolu/Vidalı"
$ws.Hyperlinks.Add($ws.Range("K10"), "https://github.com/btk42/AC-RFD-00-235-000-H3B1-01", "", "", "https://github.com/btk42/AC-RFD-00-235-000-H3B1-01") | Out-Null

# --- Old row 11 ("AC-RFD-00-235-LOP-H2B1-01") is removed entirely; rows below
#     shift up by one. ---
$ws.Rows("11:11").Delete()

# --- Re-sequence the "Sıra" numbers in column A for the rows that moved up. ---
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12

# --- Update the remembered selection to match the saved workbook state. ---
$ws.Range("N15").Select()
